# Commit: exchange mysql and dingo case order, add some hash partition dml cases
#
# Appends 9 new DML test-case rows (rows 114-122) to Sheet1, covering
# hash-partition and scalar-index-hash-partition update/delete cases for
# tables scalar056/scalar057/scalar058, and a range-partition case set for
# scalar058. Also updates the sheet selection to match the post-edit
# worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(114, 1).Value = 'updel_113'
$ws.Cells.Item(114, 2).Value = 'y'
$ws.Cells.Item(114, 3).Value = '表hash分区更新普通字段值'
$ws.Cells.Item(114, 4).Value = 'SQLFunction'
$ws.Cells.Item(114, 6).Value = 'scalar056'
$ws.Cells.Item(114, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(114, 8).Value = 'update $scalar056 set name=''Java'' where id between 5 and 20'
$ws.Cells.Item(114, 9).Value = '16'
$ws.Cells.Item(114, 10).Value = 'select id,name from $scalar056'
$ws.Cells.Item(114, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_113.csv'
$ws.Cells.Item(114, 12).Value = 'csv_containsAll'

$ws.Cells.Item(115, 1).Value = 'updel_114'
$ws.Cells.Item(115, 2).Value = 'y'
$ws.Cells.Item(115, 3).Value = '表hash分区更新索引字段值'
$ws.Cells.Item(115, 4).Value = 'Index'
$ws.Cells.Item(115, 5).Value = 'scalar_index'
$ws.Cells.Item(115, 6).Value = 'scalar056'
$ws.Cells.Item(115, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(115, 8).Value = 'update $scalar056 set age=100 where age between 10 and 100'
$ws.Cells.Item(115, 9).Value = '19'
$ws.Cells.Item(115, 10).Value = 'select id,age from $scalar056'
$ws.Cells.Item(115, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_114.csv'
$ws.Cells.Item(115, 12).Value = 'csv_containsAll'

$ws.Cells.Item(116, 1).Value = 'updel_115'
$ws.Cells.Item(116, 2).Value = 'y'
$ws.Cells.Item(116, 3).Value = '表hash分区删除数据'
$ws.Cells.Item(116, 4).Value = 'SQLFunction'
$ws.Cells.Item(116, 6).Value = 'scalar056'
$ws.Cells.Item(116, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(116, 8).Value = 'delete from $scalar056 where id in (1,3,5,7,9,10,11,13,15,17,19,21,23,25,27,29,31,33,35,37,39)'
$ws.Cells.Item(116, 9).Value = '17'
$ws.Cells.Item(116, 10).Value = 'select * from $scalar056'
$ws.Cells.Item(116, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_115.csv'
$ws.Cells.Item(116, 12).Value = 'csv_containsAll'

$ws.Cells.Item(117, 1).Value = 'updel_116'
$ws.Cells.Item(117, 2).Value = 'y'
$ws.Cells.Item(117, 3).Value = '标量索引hash分区更新普通字段值'
$ws.Cells.Item(117, 4).Value = 'Index'
$ws.Cells.Item(117, 5).Value = 'scalar_index'
$ws.Cells.Item(117, 6).Value = 'scalar057'
$ws.Cells.Item(117, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(117, 8).Value = 'update $scalar057 set name=''Java'' where id between 5 and 20'
$ws.Cells.Item(117, 9).Value = '16'
$ws.Cells.Item(117, 10).Value = 'select id,name from $scalar057'
$ws.Cells.Item(117, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_116.csv'
$ws.Cells.Item(117, 12).Value = 'csv_containsAll'

$ws.Cells.Item(118, 1).Value = 'updel_117'
$ws.Cells.Item(118, 2).Value = 'y'
$ws.Cells.Item(118, 3).Value = '标量索引hash分区更新索引字段值'
$ws.Cells.Item(118, 4).Value = 'Index'
$ws.Cells.Item(118, 5).Value = 'scalar_index'
$ws.Cells.Item(118, 6).Value = 'scalar057'
$ws.Cells.Item(118, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(118, 8).Value = 'update $scalar057 set age=100 where age between 10 and 100'
$ws.Cells.Item(118, 9).Value = '19'
$ws.Cells.Item(118, 10).Value = 'select id,age from $scalar057'
$ws.Cells.Item(118, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_117.csv'
$ws.Cells.Item(118, 12).Value = 'csv_containsAll'

$ws.Cells.Item(119, 1).Value = 'updel_118'
$ws.Cells.Item(119, 2).Value = 'y'
$ws.Cells.Item(119, 3).Value = '标量索引hash分区删除数据'
$ws.Cells.Item(119, 4).Value = 'Index'
$ws.Cells.Item(119, 5).Value = 'scalar_index'
$ws.Cells.Item(119, 6).Value = 'scalar057'
$ws.Cells.Item(119, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(119, 8).Value = 'delete from $scalar057 where id in (1,3,5,7,9,10,11,13,15,17,19,21,23,25,27,29,31,33,35,37,39)'
$ws.Cells.Item(119, 9).Value = '17'
$ws.Cells.Item(119, 10).Value = 'select * from $scalar057'
$ws.Cells.Item(119, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_118.csv'
$ws.Cells.Item(119, 12).Value = 'csv_containsAll'

$ws.Cells.Item(120, 1).Value = 'updel_119'
$ws.Cells.Item(120, 2).Value = 'y'
$ws.Cells.Item(120, 3).Value = '表range分区更新普通字段值'
$ws.Cells.Item(120, 4).Value = 'SQLFunction'
$ws.Cells.Item(120, 6).Value = 'scalar058'
$ws.Cells.Item(120, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(120, 8).Value = 'update $scalar058 set name=''Java'' where id between 5 and 20'
$ws.Cells.Item(120, 9).Value = '16'
$ws.Cells.Item(120, 10).Value = 'select id,name from $scalar058'
$ws.Cells.Item(120, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_119.csv'
$ws.Cells.Item(120, 12).Value = 'csv_containsAll'

$ws.Cells.Item(121, 1).Value = 'updel_120'
$ws.Cells.Item(121, 2).Value = 'y'
$ws.Cells.Item(121, 3).Value = '表range分区更新索引字段值'
$ws.Cells.Item(121, 4).Value = 'Index'
$ws.Cells.Item(121, 5).Value = 'scalar_index'
$ws.Cells.Item(121, 6).Value = 'scalar058'
$ws.Cells.Item(121, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(121, 8).Value = 'update $scalar058 set age=100 where age between 10 and 100'
$ws.Cells.Item(121, 9).Value = '19'
$ws.Cells.Item(121, 10).Value = 'select id,age from $scalar058'
$ws.Cells.Item(121, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_120.csv'
$ws.Cells.Item(121, 12).Value = 'csv_containsAll'

$ws.Cells.Item(122, 1).Value = 'updel_121'
$ws.Cells.Item(122, 2).Value = 'y'
$ws.Cells.Item(122, 3).Value = '表range分区删除数据'
$ws.Cells.Item(122, 4).Value = 'SQLFunction'
$ws.Cells.Item(122, 6).Value = 'scalar058'
$ws.Cells.Item(122, 7).Value = 'scalar_common_value1'
$ws.Cells.Item(122, 8).Value = 'delete from $scalar058 where id in (1,3,5,7,9,10,11,13,15,17,19,21,23,25,27,29,31,33,35,37,39)'
$ws.Cells.Item(122, 9).Value = '17'
$ws.Cells.Item(122, 10).Value = 'select * from $scalar058'
$ws.Cells.Item(122, 11).Value = 'src/test/resources/io.dingodb.test/testdata/cases/dml/updatedelete/expectedresult/updatedelete_121.csv'
$ws.Cells.Item(122, 12).Value = 'csv_containsAll'


# Move the active selection to match the post-edit workbook state.
$ws.Range("I128").Select()
